$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B10").Value = "0b011100"
$ws.Range("B11").Value = "0b010101"
$ws.Range("B12").Value = "0b010011"
